$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "160×6=960" "392×3=1176"
Replace-Text "342×9=3078" "568×7=3976"
Replace-Text "415×8=3320" "715×8=5720"
Replace-Text "856×4=3424" "367×8=2936"
Replace-Text "731×4=2924" "175×5=875"
Replace-Text "424×2=848" "163×6=978"
Replace-Text "111×4=444" "195×6=1170"
Replace-Text "863×5=4315" "558×9=5022"
Replace-Text "110×6=660" "732×2=1464"
Replace-Text "992×9=8928" "812×4=3248"
Replace-Text "245×2=490" "106×3=318"
Replace-Text "447×5=2235" "262×4=1048"
Replace-Text "980×3=2940" "906×3=2718"
Replace-Text "750×7=5250" "545×2=1090"
Replace-Text "726×8=5808" "519×6=3114"
Replace-Text "349×9=3141" "389×4=1556"
Replace-Text "118×9=1062" "190×4=760"
Replace-Text "579×2=1158" "619×3=1857"
Replace-Text "481×9=4329" "296×5=1480"
Replace-Text "875×8=7000" "491×2=982"
Replace-Text "983×2=1966" "672×4=2688"
Replace-Text "269×9=2421" "844×8=6752"
Replace-Text "423×2=846" "429×5=2145"
Replace-Text "509×4=2036" "623×9=5607"
Replace-Text "437×4=1748" "537×6=3222"
